# Applies the "Generalising plotting options and improving quality of figures" edit
# to the Setup sheet of the parameters workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# ---------------------------------------------------------------------------
# 1. Update the file-path parameters (rows 2-4) from the work-desktop paths
#    to the work-laptop paths, and the model output folder (row 5).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value2 = "M:\Working\NewModel\ModelInputs\Tar_AvMetData_1981-2010.csv"
$ws.Range("C3").Value2 = "M:\Working\NewModel\ModelInputs\obs_csvs\Coull_9amDailyMeanQ_oldRating.xlsx"
$ws.Range("C4").Value2 = "M:\Working\NewModel\ModelInputs\obs_csvs\WholePeriod\TarChem_R4_SS-P.xlsx"
$ws.Range("C5").Value2 = "M:\Working\NewModel\ModelOutputs"

# ---------------------------------------------------------------------------
# 2. Replace the old "TC_variables_to_plot" row (15) with the new, more
#    generic "plot_TC" parameter, and add three new rows (16-18) covering
#    plot_R, plot_R_reaches and R_vars_to_plot. Row 19 becomes blank again.
# ---------------------------------------------------------------------------

# Row 15: plot_TC
$ws.Range("A15").Value2 = "plot_TC"
$ws.Range("B15").Value2 = "Save a plot of simulated output from the terrestrial compartment?"
$ws.Range("C15").Value2 = "y"
$ws.Range("D15").Value2 = "Text: y or n (case sensitive)"

# Row 16: plot_R
$ws.Range("A16").Value2 = "plot_R"
$ws.Range("B16").Value2 = "Save a plot of simulated output from the stream reaches?"
$ws.Range("C16").Value2 = "y"
$ws.Range("D16").Value2 = "Text: y or n (case sensitive)"

# Row 17: plot_R_reaches
$ws.Range("A17").Value2 = "plot_R_reaches"
$ws.Range("B17").Value2 = "List of reaches to plot reach output for"
$ws.Range("C17").Value2 = "all"
$ws.Range("D17").Value2 = "Either a list (e.g. 1,2,10), or all (case-sensitive)"

# Row 18: R_vars_to_plot
$ws.Range("A18").Value2 = "R_vars_to_plot"
$ws.Range("B18").Value2 = "List of instream variables to plot"
$ws.Range("C18").Value2 = "SS, TDP, PP, TP, Q"
$ws.Range("D18").Value2 = "Choose from: SS, TDP, PP, TP, Q"

# Row 19 is now empty - clear any leftover value/format remnants.
$ws.Range("A19:D19").ClearContents()

# ---------------------------------------------------------------------------
# 3. Formatting to match the restyled block (rows 15-19).
# ---------------------------------------------------------------------------

# Row heights: rows 15 & 16 wrap onto two lines (height 30); rows 17-19 use
# the sheet's default row height again.
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 15

# Column A & D: vertical-top alignment, no wrap.
$aRange = $ws.Range("A15:A18")
$aRange.VerticalAlignment = -4160
$aRange.WrapText = $false
$aRange.HorizontalAlignment = -4142

# Column B: vertical-top, wrap text (descriptions can be long).
$bRange = $ws.Range("B15:B16")
$bRange.VerticalAlignment = -4160
$bRange.WrapText = $true
$bRange.HorizontalAlignment = -4142

$ws.Range("B17:B18").VerticalAlignment = -4160
$ws.Range("B17:B18").WrapText = $false
$ws.Range("B17:B18").HorizontalAlignment = -4142

# Column C: horizontal-left, vertical-top (values).
$cRange = $ws.Range("C15:C18")
$cRange.VerticalAlignment = -4160
$cRange.HorizontalAlignment = -4131
$cRange.WrapText = $false

# Column D (format notes) rows 15 & 16 use the normal body font/style.
$ws.Range("D15:D16").VerticalAlignment = -4160
$ws.Range("D15:D16").HorizontalAlignment = -4142
$ws.Range("D15:D16").WrapText = $false

# Column D rows 17-18 keep the distinctive "notes" font (now Calibri 11
# instead of the old Courier New 9) and are vertically centred.
$dNotes = $ws.Range("D17:D19")
$dNotes.Font.Name = "Calibri"
$dNotes.Font.Size = 11
$dNotes.VerticalAlignment = -4108

$ws.Range("A15:D19").Font.Name = "Calibri"
$ws.Range("A15:D19").Font.Size = 11

# ---------------------------------------------------------------------------
# 4. Sheet view tidy-up: selection moved to C15, freeze/topLeftCell reset.
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("C15").Select()

# ---------------------------------------------------------------------------
# 5. Page setup for the Setup sheet (now prints on A4 portrait).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6. Param_pre-processing sheet: update the selection shown when the file
#    is reopened.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Param_pre-processing")
$ws5.Application.Goto($ws5.Range("A1"), $true)
$ws5.Range("A19:A21").Select()

# Re-select the Setup sheet so it is the active tab when the file is saved,
# matching the original workbook's tabSelected state.
$ws.Select()
